$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 4 and 5 (network #2 related rows no longer exist)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# --- Row 1 (header): from_* names, net2 column dropped, gas_boiler1/2 columns added ---
$ws.Range("B1").Value = "Q_from_net1"
$ws.Range("C1").Value = "Q_from_CHP1"
$ws.Range("D1").Value = "Q_from_CHP2"
$ws.Range("E1").Value = "Q_from_solar_th1"
$ws.Range("F1").Value = "Q_from_solar_th2"
$ws.Range("G1").Value = "Q_from_pvt1"
$ws.Range("H1").Value = "Q_from_pvt2"
$ws.Range("I1").Value = "Q_from_heat_pump1"
$ws.Range("J1").Value = "Q_from_heat_pump2"
$ws.Range("K1").Value = "Q_from_gas_boiler1"

$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "Q_from_gas_boiler2"

# --- Row 2: demand1 names ---
$ws.Range("A2").Value = "param_Q_to_demand1"
$ws.Range("B2").Value = "Q_net1_demand1"
$ws.Range("C2").Value = "Q_CHP1_demand1"
$ws.Range("D2").Value = "Q_CHP2_demand1"
$ws.Range("E2").Value = "Q_solar_th1_demand1"
$ws.Range("F2").Value = "Q_solar_th2_demand1"
$ws.Range("G2").Value = "Q_pvt1_demand1"
$ws.Range("H2").Value = "Q_pvt2_demand1"
$ws.Range("I2").Value = "Q_heat_pump1_demand1"
$ws.Range("J2").Value = "Q_heat_pump2_demand1"
$ws.Range("K2").Value = "Q_gas_boiler1_demand1"
$ws.Range("L2").Value = "Q_gas_boiler2_demand1"

# --- Row 3: A3 becomes Q_to_net1, B3 becomes numeric 0, rest become *_net1 names ---
$ws.Range("A3").Value = "Q_to_net1"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "Q_CHP1_net1"
$ws.Range("D3").Value = "Q_CHP2_net1"
$ws.Range("E3").Value = "Q_solar_th1_net1"
$ws.Range("F3").Value = "Q_solar_th2_net1"
$ws.Range("G3").Value = "Q_pvt1_net1"
$ws.Range("H3").Value = "Q_pvt2_net1"
$ws.Range("I3").Value = "Q_heat_pump1_net1"
$ws.Range("J3").Value = "Q_heat_pump2_net1"
$ws.Range("K3").Value = "Q_gas_boiler1_net1"
$ws.Range("L3").Value = "Q_gas_boiler2_net1"
